# Horarios Línea 141 - scheduled scrape update (06:16:15 run)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912" (main schedule)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2, 1).Value = "Última actualización: 06:16:15"
$ws1.Cells.Item(3, 1).Value = "Total filas: 28"

# A new arrival (215A_EL PATO) was scraped before the old row 16, so insert a
# fresh row there and push the rest of the table (old rows 16-28) down by one.
$ws1.Rows.Item(16).Insert()
$ws1.Cells.Item(16, 1).Value = "06:16:15"
$ws1.Cells.Item(16, 2).Value = "06:51"
$ws1.Cells.Item(16, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(16, 4).Value = 35
$ws1.Cells.Item(16, 5).Value = "LP1912"

# Four brand-new arrivals appended at the bottom of the table.
$ws1.Cells.Item(30, 1).Value = "06:16:15"
$ws1.Cells.Item(30, 2).Value = "07:58"
$ws1.Cells.Item(30, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(30, 4).Value = 102
$ws1.Cells.Item(30, 5).Value = "LP1912"

$ws1.Cells.Item(31, 1).Value = "06:16:15"
$ws1.Cells.Item(31, 2).Value = "08:00"
$ws1.Cells.Item(31, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(31, 4).Value = 104
$ws1.Cells.Item(31, 5).Value = "LP1912"

$ws1.Cells.Item(32, 1).Value = "06:16:15"
$ws1.Cells.Item(32, 2).Value = "08:03"
$ws1.Cells.Item(32, 3).Value = "17_ROMERO"
$ws1.Cells.Item(32, 4).Value = 107
$ws1.Cells.Item(32, 5).Value = "LP1912"

$ws1.Cells.Item(33, 1).Value = "06:16:15"
$ws1.Cells.Item(33, 2).Value = "08:15"
$ws1.Cells.Item(33, 3).Value = "17_ROMERO"
$ws1.Cells.Item(33, 4).Value = 119
$ws1.Cells.Item(33, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215" (215 branch schedule)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2, 1).Value = "Última actualización: 06:16:15"
$ws2.Cells.Item(3, 1).Value = "Total filas: 4"

# Same new 215A_EL PATO arrival also belongs on the 215 branch sheet, before
# the old row 7.
$ws2.Rows.Item(7).Insert()
$ws2.Cells.Item(7, 1).Value = "06:16:15"
$ws2.Cells.Item(7, 2).Value = "06:51"
$ws2.Cells.Item(7, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(7, 4).Value = 35
$ws2.Cells.Item(7, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2, 1).Value = "Última actualización: 06:16:15"
$ws3.Cells.Item(3, 1).Value = "Total filas: 2"

# New arrival appended at the bottom.
$ws3.Cells.Item(7, 1).Value = "06:16:15"
$ws3.Cells.Item(7, 2).Value = "08:10"
$ws3.Cells.Item(7, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(7, 4).Value = 114
$ws3.Cells.Item(7, 5).Value = "L6173"
